# Updated and Added Test Cases for Multibill
# This script replays the Katalon test-run result refresh: several sheets get
# a new "Result" (Pass/Fail) and/or a new "Date" (last-run timestamp), and
# the previously active sheet (VerifyErroron2CharSearch) is swapped for the
# next one (VerifyCANSearch) as the active tab/selection.

$wb = $excel.ActiveWorkbook

# --- Result / Date updates per sheet -----------------------------------
# VerifySearchResult: Date only
$ws = $wb.Worksheets.Item("VerifySearchResult")
$ws.Range("B2").Value2 = "Mon Nov 24 16:35:36 EST 2025"

# VerifyNoModifyAmountandRedacted: Date only
$ws = $wb.Worksheets.Item("VerifyNoModifyAmountandRedacted")
$ws.Range("B2").Value2 = "Fri Nov 21 13:46:58 EST 2025"

# VerifyStaticTextOnSearch: Date only
$ws = $wb.Worksheets.Item("VerifyStaticTextOnSearch")
$ws.Range("B2").Value2 = "Mon Nov 24 16:42:59 EST 2025"

# VerifyCANSearch: Result Pass -> Fail, Date updated
$ws = $wb.Worksheets.Item("VerifyCANSearch")
$ws.Range("A2").Value2 = "Fail"
$ws.Range("B2").Value2 = "Mon Nov 24 14:36:07 EST 2025"

# VerifyStreetAddressSearch: Date only
$ws = $wb.Worksheets.Item("VerifyStreetAddressSearch")
$ws.Range("B2").Value2 = "Mon Nov 24 14:44:23 EST 2025"

# VerifyStaticTextOnViewCart: Date on both rows
$ws = $wb.Worksheets.Item("VerifyStaticTextOnViewCart")
$ws.Range("B2").Value2 = "Mon Nov 24 16:40:18 EST 2025"
$ws.Range("B3").Value2 = "Mon Nov 24 16:40:54 EST 2025"

# VerifyLookup1Search: Date only
$ws = $wb.Worksheets.Item("VerifyLookup1Search")
$ws.Range("B2").Value2 = "Mon Nov 24 14:40:13 EST 2025"

# VerifyUDF3Saerch: Date only
$ws = $wb.Worksheets.Item("VerifyUDF3Saerch")
$ws.Range("B2").Value2 = "Mon Nov 24 14:47:01 EST 2025"

# VerifyRemoveCartContent: Result Pass -> Fail, Date updated
$ws = $wb.Worksheets.Item("VerifyRemoveCartContent")
$ws.Range("A2").Value2 = "Fail"
$ws.Range("B2").Value2 = "Mon Nov 24 16:30:33 EST 2025"

# VerifyDataOnCartContent: Date only
$ws = $wb.Worksheets.Item("VerifyDataOnCartContent")
$ws.Range("B2").Value2 = "Fri Nov 21 19:53:20 EST 2025"

# VerifyPaymentEntryPageCC: Result Pass -> Fail, Date updated
$ws = $wb.Worksheets.Item("VerifyPaymentEntryPageCC")
$ws.Range("A2").Value2 = "Fail"
$ws.Range("B2").Value2 = "Mon Nov 24 15:02:44 EST 2025"

# VerifyPaymentEntryPagePC: Result Pass -> Fail, Date updated
$ws = $wb.Worksheets.Item("VerifyPaymentEntryPagePC")
$ws.Range("A2").Value2 = "Fail"
$ws.Range("B2").Value2 = "Mon Nov 24 15:18:20 EST 2025"

# VerifyPaymentEntryPageCorp: Result Pass -> Fail, Date updated
$ws = $wb.Worksheets.Item("VerifyPaymentEntryPageCorp")
$ws.Range("A2").Value2 = "Fail"
$ws.Range("B2").Value2 = "Mon Nov 24 15:14:07 EST 2025"

# --- Active tab / selection change --------------------------------------
# Previously VerifyErroron2CharSearch (tab 10, selection E2) was active;
# now VerifyCANSearch (tab 11, selection I2) is active.
$wsNext = $wb.Worksheets.Item("VerifyCANSearch")
$wsNext.Activate()
$wsNext.Range("I2").Select()
